$wb = $excel.ActiveWorkbook

# Sheet ALC - hunk at diff line 727
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 199.11111
$ws.Range("I2").Value = 206
$ws.Range("J2").Value = 175
$ws.Range("K2").Value = 206
$ws.Range("L2").Value = 175
$ws.Range("M2").Value = -93
$ws.Range("N2").Value = -401

# Sheet ALC - hunk at diff line 2028
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 436.1111
$ws.Range("I28").Value = 423.33334
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 423.33334
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 61.66665999999998
$ws.Range("N28").Value = -1470

# Sheet ALC - hunk at diff line 2640
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5910.5
$ws.Range("I40").Value = 4799.154
$ws.Range("K40").Value = 4799.154
$ws.Range("M40").Value = -4624.154

# Sheet ALC - hunk at diff line 2692
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1110.6666
$ws.Range("I41").Value = 1123.9286
$ws.Range("K41").Value = 1123.9286
$ws.Range("M41").Value = -683.9286

# Sheet ALC - hunk at diff line 2796
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7388
$ws.Range("I43").Value = 4913.1665
$ws.Range("J43").Value = 11100.25
$ws.Range("K43").Value = 4913.1665
$ws.Range("L43").Value = 11100.25
$ws.Range("M43").Value = -4844.1665
$ws.Range("N43").Value = -11238.25

# Sheet ALC - hunk at diff line 4366
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3218.5
$ws.Range("I74").Value = 3218.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3218.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2282.5
$ws.Range("N74").ClearContents()

# Sheet ALC - hunk at diff line 4516
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3218.5
$ws.Range("I77").Value = 3218.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 16092.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -11412.5
$ws.Range("N77").ClearContents()

# Sheet ALC - hunk at diff line 5523
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1512.4286
$ws.Range("J97").Value = 1512.4286
$ws.Range("L97").Value = 4537.2858
$ws.Range("N97").Value = -5529.2858

# Sheet ALC - hunk at diff line 5676
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 6750.5
$ws.Range("I100").Value = 2250
$ws.Range("J100").Value = 15751.5
$ws.Range("K100").Value = 2250
$ws.Range("L100").Value = 15751.5
$ws.Range("M100").Value = -1709
$ws.Range("N100").Value = -16833.5

# Sheet ARM - hunk at diff line 9328
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5320.727
$ws.Range("I32").Value = 4746.59
$ws.Range("K32").Value = 4746.59
$ws.Range("M32").Value = -4459.59

# Sheet ARM - hunk at diff line 13717
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2434.2307
$ws.Range("I122").Value = 1600.4615
$ws.Range("K122").Value = 4801.3845
$ws.Range("M122").Value = -2351.3845

# Sheet ARM - hunk at diff line 14201
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2540.5686
$ws.Range("I132").Value = 2064.4187
$ws.Range("K132").Value = 6193.256100000001
$ws.Range("M132").Value = -3663.256100000001

# Sheet BSM - hunk at diff line 14950
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 517.75
$ws.Range("I5").Value = 562.8
$ws.Range("J5").Value = 442.66666
$ws.Range("K5").Value = 562.8
$ws.Range("L5").Value = 442.66666
$ws.Range("M5").Value = -449.8
$ws.Range("N5").Value = -668.66666

# Sheet BSM - hunk at diff line 15682
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4427.125
$ws.Range("I20").Value = 4719.5
$ws.Range("J20").Value = 3939.8333
$ws.Range("K20").Value = 4719.5
$ws.Range("L20").Value = 3939.8333
$ws.Range("M20").Value = -4472.5
$ws.Range("N20").Value = -4433.8333

# Sheet BSM - hunk at diff line 19278
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1821.2858
$ws.Range("J94").Value = 1375
$ws.Range("L94").Value = 1375
$ws.Range("N94").Value = -2277

# Sheet BSM - hunk at diff line 20669
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 75000
$ws.Range("I123").Value = 50000
$ws.Range("J123").Value = 87500
$ws.Range("K123").Value = 50000
$ws.Range("L123").Value = 87500
$ws.Range("M123").Value = -45100
$ws.Range("N123").Value = -97300

# Sheet BSM - hunk at diff line 21542
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 72110.39999999999
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# Sheet CRP - hunk at diff line 21945
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 76923330
$ws.Range("I7").Value = 142857330
$ws.Range("K7").Value = 142857330
$ws.Range("M7").Value = -142857217

# Sheet CRP - hunk at diff line 23142
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4410.706
$ws.Range("I31").Value = 2228.3076
$ws.Range("K31").Value = 2228.3076
$ws.Range("M31").Value = -1933.3076

# Sheet CRP - hunk at diff line 23289
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4410.706
$ws.Range("I34").Value = 2228.3076
$ws.Range("K34").Value = 2228.3076
$ws.Range("M34").Value = -2026.3076

# Sheet CUL - hunk at diff line 28624
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 360.375
$ws.Range("J2").Value = 158.41667
$ws.Range("L2").Value = 950.5000200000001
$ws.Range("N2").Value = -1176.50002

# Sheet CUL - hunk at diff line 30445
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 69.75
$ws.Range("J38").Value = 127.2
$ws.Range("L38").Value = 381.6
$ws.Range("N38").Value = -1075.6

# Sheet CUL - hunk at diff line 34559
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 17749.5
$ws.Range("I120").Value = 17749.5
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 53248.5
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -48410.5
$ws.Range("N120").ClearContents()

# Sheet GSM - hunk at diff line 41647
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4582.857
$ws.Range("I122").Value = 4012.7896
$ws.Range("K122").Value = 12038.3688
$ws.Range("M122").Value = -9588.3688

# Sheet LTW - hunk at diff line 44856
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 35638.77
$ws.Range("I46").Value = 87562.2
$ws.Range("K46").Value = 87562.2
$ws.Range("M46").Value = -87374.2

# Sheet LTW - hunk at diff line 46614
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2723.3333
$ws.Range("I82").Value = 2680.3333
$ws.Range("J82").Value = 2766.3333
$ws.Range("K82").Value = 2680.3333
$ws.Range("L82").Value = 2766.3333
$ws.Range("M82").Value = -2319.3333
$ws.Range("N82").Value = -3488.3333

# Sheet LTW - hunk at diff line 46764
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2723.3333
$ws.Range("I85").Value = 2680.3333
$ws.Range("J85").Value = 2766.3333
$ws.Range("K85").Value = 2680.3333
$ws.Range("L85").Value = 2766.3333
$ws.Range("M85").Value = -1432.3333
$ws.Range("N85").Value = -5262.3333

# Sheet LTW - hunk at diff line 48556
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 204865
$ws.Range("I122").Value = 669668
$ws.Range("K122").Value = 2009004
$ws.Range("M122").Value = -2006554

# Sheet WVR - hunk at diff line 52944
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 17186.5
$ws.Range("J70").Value = 17186.5
$ws.Range("L70").Value = 17186.5
$ws.Range("N70").Value = -17816.5

# Sheet WVR - hunk at diff line 53088
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 17186.5
$ws.Range("J73").Value = 17186.5
$ws.Range("L73").Value = 17186.5
$ws.Range("N73").Value = -19370.5

# Sheet WVR - hunk at diff line 54417
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3601.0625
$ws.Range("I100").Value = 3843.25
$ws.Range("K100").Value = 7686.5
$ws.Range("M100").Value = -7145.5

# Sheet WVR - hunk at diff line 55492
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1756.909
$ws.Range("I122").Value = 1833.6
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 5500.799999999999
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = -3050.799999999999
$ws.Range("N122").Value = -7870

Write-Host "All updates applied"